$d = $word.ActiveDocument

# --- XML fragments for the three new body paragraphs (diff additions) ------------
$newPara1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Il codice IJVM </w:t></w:r><w:r><w:t xml:space="preserve">prodotto </w:t></w:r><w:r><w:t>può essere visto come il risultato di raffinamenti successivi di uno più grezzo. In una prima fase si è letto il codice C-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>like</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> linea per linea, producendone una diretta traduzione in IJVM. Già nella prima traduzione si è cercato di evitare di far uso di variabili, </w:t></w:r><w:r><w:t>a favore del</w:t></w:r><w:r><w:t xml:space="preserve">l’indirizzamento </w:t></w:r><w:r><w:t xml:space="preserve">a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. La scelta di usare lo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> è data soprattutto </w:t></w:r><w:r><w:t xml:space="preserve">da </w:t></w:r><w:r><w:t>una mera questione economica nella gestione dello spazio della memoria</w:t></w:r><w:r><w:t>, anche se non per questo le prestazioni si sarebbero degradate. Completata la fase di traduzione di un blocco di codice C-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>like</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, quindi anche solo un ciclo, verificato il suo funzionamento, con pochi e semplici test, abbiamo cercato di migliorarla. Il miglioramento consiste nella possibilità di ridurre il più possibi</w:t></w:r><w:r><w:t>le il numero di linee di codice ed</w:t></w:r><w:r><w:t xml:space="preserve"> eliminare</w:t></w:r><w:r><w:t xml:space="preserve"> le variabili utilizzate</w:t></w:r><w:r><w:t>, qualora ce ne fossero</w:t></w:r><w:r><w:t xml:space="preserve">. A volte la riduzione di codice non era strettamente necessaria, ma lo si è fatto per </w:t></w:r><w:r><w:t xml:space="preserve">lo stesso motivo introdotto per le variabili. </w:t></w:r><w:r><w:t xml:space="preserve">Il miglioramento del codice IJVM è stato possibile grazie a simulazioni (su pezzo di carta) dello </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stack</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, per verificarne la possibile fattibilità.</w:t></w:r></w:p>'
$newPara2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Si fa presente che per una più rapida implementazione del codice sorgente dei metodi che stampano stringhe a video si è fatto uso di un programma in C</w:t></w:r><w:r><w:t xml:space="preserve"> (si veda file </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>gen_code.c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>, che data in input una stringa genera il codice IJVM per stamparla</w:t></w:r><w:r><w:t xml:space="preserve">. Ogni codice ASCII in esadecimale, presente nel codice </w:t></w:r><w:r><w:t>generato</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>, riporta a lato il commento del corrispondente carattere.</w:t></w:r></w:p>'
$newPara3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Per un approfondimento riguardo alla realizzazione si invita il lettore a fare riferimento ai commenti nel file sorgente (</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>es2_GR09.jas</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p>'

# XML fragment that replaces the closing "COMMENTI CONCLUSIVI" heading paragraph,
# dropping the "_GoBack" bookmark that used to sit there (it moved up into the
# second new paragraph above, which already carries its own bookmark pair).
$closingPara = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Titolosezione"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>COMMENTI CONCLUSIVI</w:t></w:r></w:p>'

# --- locate the anchor paragraph ("REALIZZAZIONE DEL CODICE IJVM") ---------------
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "REALIZZAZIONE DEL CODICE IJVM") {
        $anchor = $p
        break
    }
}
$anchorIndex = $anchor.Index

# --- insert the three new paragraphs right after the anchor, one at a time -------
# Each is created as a bare empty paragraph first (InsertParagraphAfter) and then
# its content is replaced wholesale via InsertXML, so it never inherits the run/
# paragraph formatting of its neighbours.
[void]$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()
[void]$d.Paragraphs($anchorIndex + 1).Range.InsertXML($newPara1)

[void]$d.Paragraphs($anchorIndex + 1).Range.InsertParagraphAfter()
[void]$d.Paragraphs($anchorIndex + 2).Range.InsertXML($newPara2)

[void]$d.Paragraphs($anchorIndex + 2).Range.InsertParagraphAfter()
[void]$d.Paragraphs($anchorIndex + 3).Range.InsertXML($newPara3)

# --- drop the bookmark from the final "COMMENTI CONCLUSIVI" heading --------------
$closing = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "COMMENTI CONCLUSIVI") {
        $closing = $p
        break
    }
}
[void]$closing.Range.InsertXML($closingPara)

Write-Host "Done"
